$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete data row 2 (Trial_No 102), shifting all subsequent rows up by one
# and dropping what used to be the last row (115).
$ws.Rows.Item(2).Delete()

# The autofilter range doesn't auto-shrink on row delete in this host, so
# reapply it explicitly over the new used range (A1:K114).
$ws.AutoFilterMode = $false
$ws.Range("A1:K114").AutoFilter() | Out-Null

# Keep the workbook-level _FilterDatabase defined name in sync with the
# resized autofilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$K`$114"
    }
}

# Update the sheet's active selection to match the post-edit view state.
$ws.Range("O8").Select() | Out-Null
